$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update B2:B13 values ---
$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value = 428987.41
$metrics.Range("B3").Value = 352159.67000000004
$metrics.Range("B4").Value = 136301.85999999999
$metrics.Range("B5").Value = 17199
$metrics.Range("B6").Value = 4796118.8800000008
$metrics.Range("B7").Value = 4041978.3399999994
$metrics.Range("B8").Value = 1406903.9999999998
$metrics.Range("B9").Value = 186200
$metrics.Range("B10").Value = 33261442.680999823
$metrics.Range("B11").Value = 31317199.860000003
$metrics.Range("B12").Value = 11688612.890000004
$metrics.Range("B13").Value = 1283827

# Move the selection on the Metrics sheet to match the saved view state
$metrics.Range("D26").Select()

# --- today sheet: the B/E/F formula columns recompute automatically from
#     the Metrics! references above, so only the saved selection needs to
#     be updated explicitly ---
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("E7").Select()
